$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts every existing cell
# (values, styles, comments, etc.) one column to the right and keeps the
# picture anchored relative to its cell, exactly like the authored edit.
$ws.Columns("A").Insert()

# Populate the freshly inserted column A with the nomenclature codes for
# each storyboard row (rows 11-18).
$ws.Range("A11").Value = "c010"
$ws.Range("A12").Value = "c020"
$ws.Range("A13").Value = "c030"
$ws.Range("A14").Value = "c040"
$ws.Range("A15").Value = "c050"
$ws.Range("A16").Value = "c060"
$ws.Range("A17").Value = "c070"
$ws.Range("A18").Value = "c080"

# Reword two of the explanatory paragraphs (now in column C after the shift).
$ws.Range("C11").Value = "Les investissement lourds, les coûts de fonctionnement élevés, les durées importantes des tâches, nécessitent d'anticiper sur le pilotage au long terme." + [char]10 + "Définir une stratégie est alors nécessaire."
$ws.Range("C13").Value = "C'est un plan global qui fixe des axes directeurs." + [char]10 + "Ce sont des lignes directrices qui vont guider les décisions et les actions."

# Match the author's final view state: zoomed in, scrolled/selected on C18.
[void]$ws.Range("C18").Select()
$excel.ActiveWindow.Zoom = 186
